$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the mislabeled sub-header in B2 ("unnamed: 1_level_1" -> "total")
$ws.Range("B2").Value = "total"

# 2. Remove the now-empty "situação do domicílio" section header row (row 5).
#    Everything below shifts up by one row.
$ws.Rows("5").Delete()

# 3. Remove the now-empty "grandes regiões e unidades da federação" section
#    header row. After the previous deletion it now sits at row 7.
$ws.Rows("7").Delete()
